$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply corrected stock quantities (column F), recomputed values (column G),
# and recomputed company Sub Total / Grand Total figures (column B) exactly as
# reported in the corrected stock report.
$ws.Range("F80").Value = 22
$ws.Range("G80").Value = 6932.86
$ws.Range("B81").Value = 183445.69
$ws.Range("F122").Value = 85
$ws.Range("G122").Value = 8709.1
$ws.Range("F123").Value = 37
$ws.Range("G123").Value = 7463.27
$ws.Range("F128").Value = 3
$ws.Range("G128").Value = 807.5700000000001
$ws.Range("B133").Value = 202142.46
$ws.Range("F174").Value = 17
$ws.Range("G174").Value = 544.34
$ws.Range("B176").Value = 13204.79
$ws.Range("F190").Value = 7
$ws.Range("G190").Value = 8815.66
$ws.Range("B195").Value = 36267.13
$ws.Range("F242").Value = 13
$ws.Range("G242").Value = 624.78
$ws.Range("B243").Value = 759.5599999999999
$ws.Range("F246").Value = 73
$ws.Range("G246").Value = 5030.43
$ws.Range("B248").Value = 5561.68
$ws.Range("F295").Value = 25
$ws.Range("G295").Value = 6184.5
$ws.Range("B310").Value = 107270.13
$ws.Range("B314").Value = 61610
$ws.Range("D314").Value = 102.71
$ws.Range("E314").Value = 122.71
$ws.Range("F314").Value = 82
$ws.Range("G314").Value = 8422.219999999999
$ws.Range("B315").Value = 57077
$ws.Range("D315").Value = 93.08
$ws.Range("E315").Value = 111.2
$ws.Range("F315").Value = 1
$ws.Range("G315").Value = 93.08
$ws.Range("F327").Value = 50
$ws.Range("G327").Value = 6853.5
$ws.Range("F331").Value = 28
$ws.Range("G331").Value = 3316.88
$ws.Range("F344").Value = 121
$ws.Range("G344").Value = 15610.21
$ws.Range("F354").Value = 25
$ws.Range("G354").Value = 2527.5
$ws.Range("F370").Value = 10
$ws.Range("G370").Value = 1997.7
$ws.Range("F372").Value = 74
$ws.Range("G372").Value = 5202.94
$ws.Range("B380").Value = 239334.79
$ws.Range("F385").Value = 5
$ws.Range("G385").Value = 1016.55
$ws.Range("B389").Value = 21063.58
$ws.Range("F435").Value = 73
$ws.Range("G435").Value = 4009.16
$ws.Range("F436").Value = 61
$ws.Range("G436").Value = 6437.94
$ws.Range("F443").Value = 25
$ws.Range("G443").Value = 930.25
$ws.Range("B447").Value = 37548.68
$ws.Range("F490").Value = 452
$ws.Range("G490").Value = 6079.4
$ws.Range("F491").Value = 460
$ws.Range("G491").Value = 6049
$ws.Range("F492").Value = 540
$ws.Range("G492").Value = 6917.4
$ws.Range("F495").Value = 297
$ws.Range("G495").Value = 4879.71
$ws.Range("F501").Value = 912
$ws.Range("G501").Value = 6000.96
$ws.Range("F503").Value = 925
$ws.Range("G503").Value = 6003.25
$ws.Range("F507").Value = 537
$ws.Range("G507").Value = 7910.01
$ws.Range("B508").Value = 92551.16
$ws.Range("F510").Value = 9
$ws.Range("G510").Value = 332.73
$ws.Range("F513").Value = 31
$ws.Range("G513").Value = 1079.11
$ws.Range("B515").Value = 6108.11
$ws.Range("F539").Value = 6
$ws.Range("G539").Value = 2510.16
$ws.Range("F540").Value = 4
$ws.Range("G540").Value = 2006.8
$ws.Range("F542").Value = 5
$ws.Range("G542").Value = 3313.85
$ws.Range("F544").Value = 20
$ws.Range("G544").Value = 3287.2
$ws.Range("B546").Value = 13295.74
$ws.Range("F555").Value = 348
$ws.Range("G555").Value = 2383.8
$ws.Range("F557").Value = 585
$ws.Range("G557").Value = 11612.25
$ws.Range("F558").Value = 234
$ws.Range("G558").Value = 1567.8
$ws.Range("B562").Value = 35487.41
$ws.Range("F574").Value = 13
$ws.Range("G574").Value = 678.34
$ws.Range("B583").Value = 22508.96
$ws.Range("F619").Value = 34
$ws.Range("G619").Value = 3486.36
$ws.Range("B639").Value = 205725.07
$ws.Range("F672").Value = 47
$ws.Range("G672").Value = 2029.46
$ws.Range("F673").Value = 10
$ws.Range("G673").Value = 331.1
$ws.Range("B676").Value = 19163.74
$ws.Range("F695").Value = 9
$ws.Range("G695").Value = 557.1
$ws.Range("F705").Value = 3
$ws.Range("G705").Value = 2840.13
$ws.Range("B713").Value = 77811.17
$ws.Range("F758").Value = 22
$ws.Range("G758").Value = 2453.88
$ws.Range("B772").Value = 227470.14
$ws.Range("F819").Value = 0
$ws.Range("G819").Value = 0
$ws.Range("F826").Value = 5
$ws.Range("G826").Value = 7177.05
$ws.Range("B840").Value = 86658.99000000001
$ws.Range("F849").Value = 509
$ws.Range("G849").Value = 15387.07
$ws.Range("F857").Value = 96
$ws.Range("G857").Value = 12343.68
$ws.Range("B858").Value = 598045.84
$ws.Range("B864").Value = 3312541.93
$ws.Range("B865").Value = 3312541.93
